$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21646
$ws.Range("E2").Value = 2188
$ws.Range("J2").Value = 5067

$ws.Range("B3").Value = 32992
$ws.Range("E3").Value = 2429
$ws.Range("J3").Value = 6046

$ws.Range("H4").Value = 18600
$ws.Range("J4").Value = 6987

$ws.Range("F5").Value = 38960
$ws.Range("J5").Value = 7742

$ws.Range("D6").Value = 22815
$ws.Range("J6").Value = 8997

$ws.Range("J7").Value = 10674

$ws.Range("D8").Value = 19906
$ws.Range("F8").Value = 42721
$ws.Range("J8").Value = 8404

$ws.Range("F9").Value = 58830
$ws.Range("J9").Value = 10821

$ws.Range("B10").Value = 81587
$ws.Range("C10").Value = 49203
$ws.Range("J10").Value = 12899

$ws.Range("B11").Value = 77773
$ws.Range("C11").Value = 45979
$ws.Range("J11").Value = 12322

$ws.Range("B12").Value = 76812
$ws.Range("C12").Value = 43704
$ws.Range("E12").Value = 5652
$ws.Range("F12").Value = 79542
$ws.Range("J12").Value = 12113

$ws.Range("B13").Value = 75205
$ws.Range("C13").Value = 40560
$ws.Range("D13").Value = 29020
$ws.Range("E13").Value = 5625
$ws.Range("F13").Value = 72965
$ws.Range("H13").Value = 38849
$ws.Range("J13").Value = 10693

$ws.Range("B14").Value = 62227
$ws.Range("C14").Value = 32480
$ws.Range("D14").Value = 24531
$ws.Range("E14").Value = 5215
$ws.Range("G14").Value = 17860
$ws.Range("J14").Value = 9520

$ws.Range("B15").Value = 60788
$ws.Range("C15").Value = 30583
$ws.Range("D15").Value = 24132
$ws.Range("E15").Value = 6073
$ws.Range("G15").Value = 17776
$ws.Range("H15").Value = 28262
$ws.Range("J15").Value = 9528

$ws.Range("B16").Value = 68958
$ws.Range("C16").Value = 37002
$ws.Range("D16").Value = 25935
$ws.Range("E16").Value = 6021
$ws.Range("F16").Value = 65277
$ws.Range("H16").Value = 30922
$ws.Range("J16").Value = 9783

$ws.Range("B17").Value = 74773
$ws.Range("C17").Value = 39191
$ws.Range("E17").Value = 6481
$ws.Range("J17").Value = 9933

$ws.Range("B18").Value = 68892
$ws.Range("C18").Value = 35400
$ws.Range("D18").Value = 26671
$ws.Range("E18").Value = 6820
$ws.Range("G18").Value = 20157
$ws.Range("I18").Value = 15179
$ws.Range("J18").Value = 9602

$ws.Range("B19").Value = 73351
$ws.Range("C19").Value = 41530
$ws.Range("E19").Value = 6445
$ws.Range("F19").Value = 59243
$ws.Range("H19").Value = 29852
$ws.Range("J19").Value = 6623

$wb.Save()
